{"js": "// Replace each three-digit-division answer cell's text with its new value.\n// Old/new pairs are all unique and don't overlap, so simple sequential\n// search+replace (order independent) is safe.\nconst replacements = [\n  [\"517\u00f77=73, 6\", \"605\u00f74=151, 1\"],\n  [\"844\u00f77=120, 4\", \"178\u00f73=59, 1\"],\n  [\"860\u00f75=172, 0\", \"759\u00f77=108, 3\"],\n  [\"615\u00f75=123, 0\", \"751\u00f79=83, 4\"],\n  [\"927\u00f78=115, 7\", \"198\u00f72=99, 0\"],\n  [\"649\u00f72=324, 1\", \"234\u00f79=26, 0\"],\n  [\"153\u00f72=76, 1\", \"220\u00f76=36, 4\"],\n  [\"331\u00f77=47, 2\", \"355\u00f74=88, 3\"],\n  [\"446\u00f72=223, 0\", \"842\u00f72=421, 0\"],\n  [\"464\u00f78=58, 0\", \"514\u00f72=257, 0\"],\n  [\"951\u00f73=317, 0\", \"728\u00f79=80, 8\"],\n  [\"745\u00f78=93, 1\", \"248\u00f78=31, 0\"],\n  [\"207\u00f72=103, 1\", \"986\u00f72=493, 0\"],\n  [\"445\u00f75=89, 0\", \"250\u00f72=125, 0\"],\n  [\"126\u00f77=18, 0\", \"768\u00f78=96, 0\"],\n  [\"654\u00f78=81, 6\", \"115\u00f77=16, 3\"],\n  [\"762\u00f76=127, 0\", \"129\u00f72=64, 1\"],\n  [\"393\u00f72=196, 1\", \"321\u00f72=160, 1\"],\n  [\"394\u00f76=65, 4\", \"398\u00f72=199, 0\"],\n  [\"510\u00f77=72, 6\", \"569\u00f74=142, 1\"],\n  [\"688\u00f72=344, 0\", \"588\u00f72=294, 0\"],\n  [\"782\u00f74=195, 2\", \"758\u00f72=379, 0\"],\n  [\"585\u00f78=73, 1\", \"772\u00f76=128, 4\"],\n  [\"199\u00f79=22, 1\", \"969\u00f73=323, 0\"],\n  [\"407\u00f77=58, 1\", \"174\u00f77=24, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division answer cell's text with its new value.\n# Old/new pairs are all unique and don't overlap, so simple sequential\n# Find/Replace (order independent) is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"517\u00f77=73, 6\", \"605\u00f74=151, 1\"),\n    @(\"844\u00f77=120, 4\", \"178\u00f73=59, 1\"),\n    @(\"860\u00f75=172, 0\", \"759\u00f77=108, 3\"),\n    @(\"615\u00f75=123, 0\", \"751\u00f79=83, 4\"),\n    @(\"927\u00f78=115, 7\", \"198\u00f72=99, 0\"),\n    @(\"649\u00f72=324, 1\", \"234\u00f79=26, 0\"),\n    @(\"153\u00f72=76, 1\", \"220\u00f76=36, 4\"),\n    @(\"331\u00f77=47, 2\", \"355\u00f74=88, 3\"),\n    @(\"446\u00f72=223, 0\", \"842\u00f72=421, 0\"),\n    @(\"464\u00f78=58, 0\", \"514\u00f72=257, 0\"),\n    @(\"951\u00f73=317, 0\", \"728\u00f79=80, 8\"),\n    @(\"745\u00f78=93, 1\", \"248\u00f78=31, 0\"),\n    @(\"207\u00f72=103, 1\", \"986\u00f72=493, 0\"),\n    @(\"445\u00f75=89, 0\", \"250\u00f72=125, 0\"),\n    @(\"126\u00f77=18, 0\", \"768\u00f78=96, 0\"),\n    @(\"654\u00f78=81, 6\", \"115\u00f77=16, 3\"),\n    @(\"762\u00f76=127, 0\", \"129\u00f72=64, 1\"),\n    @(\"393\u00f72=196, 1\", \"321\u00f72=160, 1\"),\n    @(\"394\u00f76=65, 4\", \"398\u00f72=199, 0\"),\n    @(\"510\u00f77=72, 6\", \"569\u00f74=142, 1\"),\n    @(\"688\u00f72=344, 0\", \"588\u00f72=294, 0\"),\n    @(\"782\u00f74=195, 2\", \"758\u00f72=379, 0\"),\n    @(\"585\u00f78=73, 1\", \"772\u00f76=128, 4\"),\n    @(\"199\u00f79=22, 1\", \"969\u00f73=323, 0\"),\n    @(\"407\u00f77=58, 1\", \"174\u00f77=24, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
